# Apply the InputSchedule changes: introduce carbon-motivated case, shift
# the Electric Heater / Agile Tariff / Fourth cases along, and flag the
# Datum + Carbon-Motivated cases as "Smart Home".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Datum, case 1) ---
$ws.Range("AJ2").Value = 1

# --- Row 3 (Datum, case 1) ---
$ws.Range("Q3").Value = "AgileExtract2.csv"
$ws.Range("AJ3").Value = 1

# --- Row 4 (was "Electric Heater" -> now "Carbon Motivated", case 2) ---
$ws.Range("B4").Value = "Carbon Motivated"
$ws.Range("E4").Value = 0
$ws.Range("Q4").Value = "AgileExtract2.csv"
$ws.Range("AD4").Value = "EV"
$ws.Range("AE4").Value = "Carbon"
$ws.Range("AG4").Value = "Gas"
$ws.Range("AJ4").Value = 1

# --- Row 5 (was "Electric Heater" -> now "Carbon Motivated", case 2) ---
$ws.Range("B5").Value = "Carbon Motivated"
$ws.Range("E5").Value = 0
$ws.Range("Q5").Value = "AgileExtract2.csv"
$ws.Range("AE5").Value = "Carbon"
$ws.Range("AG5").Value = "Gas"
$ws.Range("AJ5").Value = 1

# --- Row 6 (was "Agile Tariff" -> now "Electric Heater", case 3) ---
$ws.Range("B6").Value = "Electric Heater"
$ws.Range("E6").Value = 10000
$ws.Range("Q6").Value = "Fixed22Tariff.csv"
$ws.Range("AG6").Value = "Electric"

# --- Row 7 (was "Agile Tariff" -> now "Electric Heater", case 3) ---
$ws.Range("B7").Value = "Electric Heater"
$ws.Range("E7").Value = 10000
$ws.Range("Q7").Value = "Fixed22Tariff.csv"
$ws.Range("AG7").Value = "Electric"

# --- Row 8 (was "Fourth" case 3 -> now "Agile Tariff" case 4) ---
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "Agile Tariff"

# --- Row 9 (was "Fourth" case 3 -> now "Agile Tariff" case 4) ---
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Agile Tariff"

# Restore the selected cell as it was left in the saved workbook
$ws.Range("AI3").Select()
